$wb = $excel.ActiveWorkbook
$wsDBD = $wb.Worksheets.Item("DBD")

# Change the field type for CreateDate (row 12) and LastUpdate (row 14)
# from DATE to TIMESTAMP
$wsDBD.Range("D12").Value = "TIMESTAMP"
$wsDBD.Range("D14").Value = "TIMESTAMP"

# Make DBD the active sheet/tab and move the selection to D14
$wsDBD.Activate() | Out-Null
$wsDBD.Range("D14").Select() | Out-Null
